# Update dSF column (column F) values for several rows, per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -3
$ws.Range("F6").Value = -4
$ws.Range("F11").Value = -10
$ws.Range("F17").Value = -1
